$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 727, shifting rows 727:768 down to 728:769
$ws.Rows(727).Insert()

# Populate the newly inserted row with the new data entry.
# Column A holds a date-like label that must stay plain text (not get
# auto-converted into a date serial number), so force text format first
# and restore the default "Normal" style afterwards so it matches the
# other (unstyled) data rows.
$ws.Range("A727").NumberFormat = "@"
$ws.Range("A727").Value = "2026/01/28"
$ws.Range("A727").Style = "Normal"

$ws.Range("B727").Value = "水"
$ws.Range("C727").Value = 3
$ws.Range("D727").Value = 201
